$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume(1h) (E) columns remain plain text so that
# values like "1.00" are not coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.671.09"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.701.08"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "682.69"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "160.73"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.496"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "0.440"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "4.325.79"
$ws.Range("D14").Value = "32.50"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "3.718.34"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "69.571.76"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "16.07"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "6.46"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "471.63"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "9.93"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "0.653"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "80.33"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "3.848.87"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "0.0000125"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "10.99"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").Value = "1.75"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "2.01"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").Value = "6.58"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "27.01"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").Value = "3.691.13"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").Value = "0.162"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").Value = "8.37"
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("D38").Value = "6.31"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("D39").Value = "2.31"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "169.43"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").Value = "0.944"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "47.29"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "28.98"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").Value = "2.74"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "0.000281"
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("D49").Value = "1.11"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "1.30"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "7.87"
$ws.Range("E51").Value = "  -0.16%  "